# feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" worksheet (fund-holdings detail, same shape as the
# other quarterly sheets) right before the "总计" (totals) sheet, and add a
# corresponding "2022-Q1" row to the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right before "总计".
# ---------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheetBefore)
$ws.Name = "2022-Q1"

# NOTE: sheet handles above are position-anchored, so after inserting the
# new sheet the old "$totalSheetBefore" handle now refers to the newly
# inserted sheet, not "总计" any more. Re-fetch "总计" by name so the rest
# of the script operates on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (row 1) - bold/centered/bordered header style, matching the
# other quarterly sheets' header row formatting.
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (row 2..21) - fund code / name / scale / position figures are
# stored as text (leading zeros, fixed decimals) just like in the sibling
# quarterly sheets; the index column (A) and rank column (H) are numbers.
$fundRows = @(
    @("012748", "华泰柏瑞远见智选混合型证券投资基金A", "26.44", "94.15", "3.49", "0.9228", 10),
    @("213003", "宝盈策略增长混合", "10.28", "94.38", "6.86", "0.7052", 4),
    @("460001", "华泰柏瑞盛世中国混合", "18.38", "85.97", "2.89", "0.5312", 10),
    @("007306", "华泰柏瑞基本面智选混合A", "3.54", "93.61", "4.35", "0.1540", 9),
    @("012749", "华泰柏瑞远见智选混合型证券投资基金C", "3.08", "94.15", "3.49", "0.1075", 10),
    @("006511", "博道卓远混合A", "2.20", "82.79", "2.94", "0.0647", 9),
    @("580006", "东吴新经济混合", "1.17", "91.56", "4.80", "0.0562", 6),
    @("007825", "博道志远混合A", "1.67", "82.63", "3.02", "0.0504", 9),
    @("561910", "招商中证电池主题交易型开放式指数证券投资基金", "1.92", "98.76", "2.50", "0.0480", 10),
    @("007307", "华泰柏瑞基本面智选混合C", "1.03", "93.61", "4.35", "0.0448", 9),
    @("005090", "嘉合睿金定期开放灵活配置混合型发起式A", "1.25", "52.49", "2.26", "0.0282", 5),
    @("007826", "博道志远混合C", "0.73", "82.63", "3.02", "0.0220", 9),
    @("005091", "嘉合睿金定期开放灵活配置混合型发起式C", "0.80", "52.49", "2.26", "0.0181", 5),
    @("003835", "鹏华沪深港新兴成长灵活配置混合", "0.61", "82.70", "2.48", "0.0151", 10),
    @("159918", "嘉实中创400ETF", "0.75", "99.13", "1.01", "0.0076", 2),
    @("006009", "国融融银灵活配置混合A", "0.05", "94.54", "6.78", "0.0034", 8),
    @("006512", "博道卓远混合C", "0.10", "82.79", "2.94", "0.0029", 9),
    @("006010", "国融融银灵活配置混合C", "0.03", "94.54", "6.78", "0.0020", 8),
    @("000926", "中信建投睿信灵活配置混合A", "0.13", "40.35", "1.45", "0.0019", 10),
    @("004676", "中信建投睿信灵活配置混合C", "0.02", "40.35", "1.45", "0.0003", 10)
)

# Pre-format the text columns as "@" so numeric-looking strings (fund codes
# with leading zeros, decimal figures) are kept as text, matching the
# original workbook's inlineStr cells instead of being coerced to numbers.
$ws.Range("B2:G21").NumberFormat = "@"

# Index column (A) - same bold/centered/bordered look as the header.
$indexRange = $ws.Range("A2:A21")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

$r = 2
foreach ($row in $fundRows) {
    $ws.Cells.Item($r, 1).NumberFormat = "General"
    $ws.Cells.Item($r, 1).Value = ($r - 2)
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).NumberFormat = "General"
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Add the "2022-Q1" summary row at the top of the "总计" data (row 2),
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$lastRow = $totalSheet.Cells.Item($totalSheet.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

# Only insert a single (1-row-tall) range so existing rows shift down by
# exactly one row, not by the height of the whole existing block.
$totalSheet.Range("A2:D2").Insert(-4121)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 20
$totalSheet.Cells.Item(2, 4).Value = 2.79

# Renumber the index column (A) for the rows that got shifted down, so it
# keeps counting 0,1,2,3,4,5 top to bottom.
$newLastRow = $lastRow + 1
for ($rr = 3; $rr -le $newLastRow; $rr++) {
    $totalSheet.Cells.Item($rr, 1).Value = $rr - 2
}
